$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 571015
$ws.Range("J17").Value = 571015
$ws.Range("L17").Value = 1713045
$ws.Range("N17").Value = -1713381

$ws.Range("H32").Value = 4875.5
$ws.Range("J32").Value = 4875.5
$ws.Range("L32").Value = 4875.5
$ws.Range("N32").Value = -5527.5

$ws.Range("H106").Value = 4334.1665
$ws.Range("I106").Value = 4001.25
$ws.Range("J106").Value = 5000
$ws.Range("K106").Value = 4001.25
$ws.Range("L106").Value = 5000
$ws.Range("M106").Value = -3370.25
$ws.Range("N106").Value = -6262

$ws.Range("H138").Value = 1572.28
$ws.Range("I138").Value = 897.89655
$ws.Range("J138").Value = 1847.7324
$ws.Range("K138").Value = 2693.68965
$ws.Range("L138").Value = 5543.197200000001
$ws.Range("M138").Value = 2446.31035
$ws.Range("N138").Value = -15823.1972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8850.407999999999
$ws.Range("I32").Value = 5457.5303
$ws.Range("J32").Value = 37011.3
$ws.Range("K32").Value = 5457.5303
$ws.Range("L32").Value = 37011.3
$ws.Range("M32").Value = -5170.5303
$ws.Range("N32").Value = -37585.3

$ws.Range("H74").Value = 1102.2433
$ws.Range("I74").Value = 1059.8387
$ws.Range("J74").Value = 1321.3334
$ws.Range("K74").Value = 1059.8387
$ws.Range("L74").Value = 1321.3334
$ws.Range("M74").Value = -185.8387
$ws.Range("N74").Value = -3069.3334

$ws.Range("H77").Value = 1102.2433
$ws.Range("I77").Value = 1059.8387
$ws.Range("J77").Value = 1321.3334
$ws.Range("K77").Value = 5299.1935
$ws.Range("L77").Value = 6606.666999999999
$ws.Range("M77").Value = -931.1935000000003
$ws.Range("N77").Value = -15342.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 28500
$ws.Range("J112").Value = 28500
$ws.Range("L112").Value = 28500
$ws.Range("N112").Value = -31454

$ws.Range("H134").Value = 115429.664
$ws.Range("I134").Value = 4648.933
$ws.Range("K134").Value = 13946.799
$ws.Range("M134").Value = -11411.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2254.1404
$ws.Range("I31").Value = 2528.162
$ws.Range("K31").Value = 2528.162
$ws.Range("M31").Value = -2233.162

$ws.Range("H34").Value = 2254.1404
$ws.Range("I34").Value = 2528.162
$ws.Range("K34").Value = 2528.162
$ws.Range("M34").Value = -2326.162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3380
$ws.Range("I56").Value = 3380
$ws.Range("K56").Value = 3380
$ws.Range("M56").Value = -2850

$ws.Range("H87").Value = 11674.75
$ws.Range("I87").Value = 5949.5
$ws.Range("J87").Value = 17400
$ws.Range("K87").Value = 17848.5
$ws.Range("L87").Value = 52200
$ws.Range("M87").Value = -16600.5
$ws.Range("N87").Value = -54696

$ws.Range("H90").Value = 11674.75
$ws.Range("I90").Value = 5949.5
$ws.Range("J90").Value = 17400
$ws.Range("K90").Value = 53545.5
$ws.Range("L90").Value = 156600
$ws.Range("M90").Value = -47305.5
$ws.Range("N90").Value = -169080

$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 6000
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -8350

$ws.Range("H122").Value = 647.8919
$ws.Range("I122").Value = 366.5
$ws.Range("J122").Value = 702.35486
$ws.Range("K122").Value = 3298.5
$ws.Range("L122").Value = 6321.193740000001
$ws.Range("M122").Value = -848.5
$ws.Range("N122").Value = -11221.19374

$ws.Range("H131").Value = 8214012
$ws.Range("I131").Value = 83500520
$ws.Range("J131").Value = 938.34546
$ws.Range("K131").Value = 250501560
$ws.Range("L131").Value = 2815.03638
$ws.Range("M131").Value = -250496520
$ws.Range("N131").Value = -12895.03638

$ws.Range("H137").Value = 44209.16
$ws.Range("J137").Value = 71557.92999999999
$ws.Range("L137").Value = 214673.79
$ws.Range("N137").Value = -224873.79

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4268.8945
$ws.Range("I70").Value = 4520.4
$ws.Range("J70").Value = 3989.4443
$ws.Range("K70").Value = 4520.4
$ws.Range("L70").Value = 3989.4443
$ws.Range("M70").Value = -4250.4
$ws.Range("N70").Value = -4529.4443

$ws.Range("H73").Value = 4268.8945
$ws.Range("I73").Value = 4520.4
$ws.Range("J73").Value = 3989.4443
$ws.Range("K73").Value = 4520.4
$ws.Range("L73").Value = 3989.4443
$ws.Range("M73").Value = -3584.4
$ws.Range("N73").Value = -5861.4443

$ws.Range("H97").Value = 1623.7778
$ws.Range("I97").Value = 1601.6471
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1601.6471
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1105.6471
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 32000
$ws.Range("J106").Value = 32000
$ws.Range("L106").Value = 32000
$ws.Range("N106").Value = -34524

$ws.Range("H132").Value = 2707.5
$ws.Range("I132").Value = 2231.0527
$ws.Range("J132").Value = 4000.7144
$ws.Range("K132").Value = 6693.158100000001
$ws.Range("L132").Value = 12002.1432
$ws.Range("M132").Value = -4163.158100000001
$ws.Range("N132").Value = -17062.1432

$ws.Range("H136").Value = 4397.524
$ws.Range("I136").Value = 2388.6924
$ws.Range("J136").Value = 7661.875
$ws.Range("K136").Value = 7166.0772
$ws.Range("L136").Value = 22985.625
$ws.Range("M136").Value = -4616.0772
$ws.Range("N136").Value = -28085.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H93").Value = 22346.334
$ws.Range("J93").Value = 22346.334
$ws.Range("L93").Value = 22346.334
$ws.Range("N93").Value = -27338.334

$ws.Range("H96").Value = 3062.7273
$ws.Range("I96").Value = 3154.4443
$ws.Range("J96").Value = 2650
$ws.Range("K96").Value = 3154.4443
$ws.Range("L96").Value = 2650
$ws.Range("M96").Value = -1781.4443
$ws.Range("N96").Value = -5396
